$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 828
$ws.Range("F3").Value = 14769
$ws.Range("F5").Value = 1683
$ws.Range("F7").Value = 2141
$ws.Range("F8").Value = 1322
$ws.Range("F9").Value = 2007
$ws.Range("F10").Value = 954
$ws.Range("F12").Value = 2371
$ws.Range("F13").Value = 629
$ws.Range("F15").Value = 3727
$ws.Range("F17").Value = 353
$ws.Range("F18").Value = 2770
$ws.Range("F19").Value = 723
$ws.Range("F22").Value = 16
$ws.Range("F23").Value = 1957
$ws.Range("F24").Value = 1148
$ws.Range("F25").Value = 1694
$ws.Range("F26").Value = 352
$ws.Range("F27").Value = 186
$ws.Range("F28").Value = 7783
$ws.Range("F29").Value = 5365
$ws.Range("F30").Value = 338
$ws.Range("F32").Value = 734
$ws.Range("F34").Value = 3437
$ws.Range("F37").Value = 368
$ws.Range("F38").Value = 164
$ws.Range("F39").Value = 139
$ws.Range("F40").Value = 4549
$ws.Range("F41").Value = 771
$ws.Range("F42").Value = 42
$ws.Range("F43").Value = 363
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 84
$ws.Range("F16").Value = 23
$ws.Range("F18").Value = 138
$ws.Range("F23").Value = 71
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8109
$ws.Range("F3").Value = 339
$ws.Range("F4").Value = 1186
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8109
$ws.Range("F3").Value = 828
$ws.Range("F4").Value = 339
$ws.Range("F5").Value = 1186
$ws.Range("F6").Value = 14769
$ws.Range("F8").Value = 84
$ws.Range("F9").Value = 1683
$ws.Range("F11").Value = 1322
$ws.Range("F12").Value = 2007
$ws.Range("F13").Value = 954
$ws.Range("F15").Value = 629
$ws.Range("F17").Value = 3727
$ws.Range("F18").Value = 353
$ws.Range("F19").Value = 2770
$ws.Range("F20").Value = 723
$ws.Range("F22").Value = 1957
$ws.Range("F27").Value = 1148
$ws.Range("F29").Value = 1694
$ws.Range("F30").Value = 352
$ws.Range("F31").Value = 186
$ws.Range("F32").Value = 7783
$ws.Range("F33").Value = 5365
$ws.Range("F35").Value = 338
$ws.Range("F36").Value = 734
$ws.Range("F38").Value = 3437
$ws.Range("F40").Value = 368
$ws.Range("F41").Value = 164
$ws.Range("F43").Value = 139
$ws.Range("F44").Value = 4549
$ws.Range("F45").Value = 771
$ws.Range("F46").Value = 42
$ws.Range("F47").Value = 363
